# Alert for client users: added case for both cases
# (1. checkbox to send to all users and 2. send a specific user only)
#
# Appends newly captured submission-time samples to the bottom of each of
# the four tracking sheets (mirrors the rows a subsequent test run would
# have produced).

$wb = $excel.ActiveWorkbook

function Add-Row {
    param(
        $ws,
        [int]$row,
        $a, $b, $c, $d, $e
    )
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
}

# ---------------------------------------------------------------------
# Sheet "Submit orders": rows 74-76
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Submit orders")

Add-Row $ws1 74 "10.03.2022 09:04 (Kyiv+Israel) 06:04 (UTC) 15:04 (Japan) 11:34 (India)" "***" "***" 1.311 -0.173
Add-Row $ws1 75 "10.03.2022 09:22 (Kyiv+Israel) 06:22 (UTC) 15:22 (Japan) 11:52 (India)" 1.467 -0.6540000000000001 "***" "***"
Add-Row $ws1 76 "10.04.2022 20:15 (Kyiv+Israel) 17:15 (UTC) 02:15 (Japan) 22:45 (India)" 1.276 -0.4630000000000001 "***" "***"

# ---------------------------------------------------------------------
# Sheet "Submit internet survey": row 71
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Submit internet survey")

Add-Row $ws2 71 "10.03.2022 09:07 (Kyiv+Israel) 06:07 (UTC) 15:07 (Japan) 11:37 (India)" "***" "***" 0.9379999999999999 -0.1729999999999999

# ---------------------------------------------------------------------
# Sheet "Submit a phone survey": row 64
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Submit a phone survey")

Add-Row $ws3 64 "10.03.2022 09:10 (Kyiv+Israel) 06:10 (UTC) 15:10 (Japan) 11:40 (India)" "***" "***" 2.056 -0.252

# ---------------------------------------------------------------------
# Sheet "Checkertificate": rows 66-74
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Checkertificate")

Add-Row $ws4 66 "10.03.2022 09:13 (Kyiv+Israel) 06:13 (UTC) 15:13 (Japan) 11:43 (India)" "***" "***" 1.039 -0.1169999999999999
Add-Row $ws4 67 "10.04.2022 11:21 (Kyiv+Israel) 08:21 (UTC) 17:21 (Japan) 13:51 (India)" 0.843 -0.117 "***" "***"
Add-Row $ws4 68 "10.04.2022 11:37 (Kyiv+Israel) 08:37 (UTC) 17:37 (Japan) 14:07 (India)" 1.019 -0.2929999999999999 "***" "***"
Add-Row $ws4 69 "10.04.2022 11:44 (Kyiv+Israel) 08:44 (UTC) 17:44 (Japan) 14:14 (India)" 0.694 0.03200000000000003 "***" "***"
Add-Row $ws4 70 "10.04.2022 11:55 (Kyiv+Israel) 08:55 (UTC) 17:55 (Japan) 14:25 (India)" 0.778 -0.08400000000000007 "***" "***"
Add-Row $ws4 71 "10.04.2022 11:57 (Kyiv+Israel) 08:57 (UTC) 17:57 (Japan) 14:27 (India)" "***" "***" 1.049 -0.1269999999999999
Add-Row $ws4 72 "10.04.2022 12:05 (Kyiv+Israel) 09:05 (UTC) 18:05 (Japan) 14:35 (India)" 0.842 -0.148 "***" "***"
Add-Row $ws4 73 "10.04.2022 12:08 (Kyiv+Israel) 09:08 (UTC) 18:08 (Japan) 14:38 (India)" "***" "***" 0.959 -0.03699999999999992
Add-Row $ws4 74 "10.04.2022 14:04 (Kyiv+Israel) 11:04 (UTC) 20:04 (Japan) 16:34 (India)" 0.725 -0.03100000000000003 "***" "***"
